# LONG XUYEN 8 - 2024 report update
# - Add new revenue entries (HD-LUXURY 622 "Phun môi" and HD-LUXURY 634) to
#   "CHI TIẾT DOANH THU" and refresh its totals row.
# - Add a new debt-collection entry to "CHI TIẾT VỀ THU NỢ" and refresh totals.
# - Add a "Ghi chú" column and 3 new expense entries to "CHI TIẾT CHI TIÊU",
#   refresh totals.
# - Refresh "DOANH SỐ CÁ NHÂN" per-employee rows (new employee Đặng Ngọc Mai,
#   drop the now-unused extra totals row).
# - Refresh "CHI TIÊU TỔNG HỢP" category roll-up (new "Ứng Lương" category).
# - Refresh "LŨY KẾ NGÀY" daily roll-up with the new days of activity.
# - Refresh "QUỸ LƯƠNG" payroll totals per employee.
# - Refresh "LỢI NHUẬN" profit summary.

$wb = $excel.ActiveWorkbook

# Helper: write a value that must stay literal TEXT even though it looks like
# a date/number (Excel would otherwise auto-convert "08-04-2024" etc. into a
# date serial). Prefixing with an apostrophe forces text, matching how the
# existing date-like cells in this report are already stored.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
}

# ============================================================
# Sheet 1: CHI TIẾT DOANH THU
# ============================================================
$ws1 = $wb.Worksheets.Item("CHI TIẾT DOANH THU")

# Row 4 used to be the "Tổng" row; it becomes a real data row.
Set-TextValue $ws1.Range("A4") "08-04-2024"
$ws1.Range("B4").Value = "HD-LUXURY"
$ws1.Range("C4").Value = 622
$ws1.Range("D4").Value = "LONG XUYÊN"
$ws1.Range("E4").Value = "Phun môi"
$ws1.Range("F4").Value = "Chị vui"
$ws1.Range("G4").Value = "CTV"
$ws1.Range("H4").Value = $null
$ws1.Range("I4").Value = 2000000
$ws1.Range("J4").Value = $null
$ws1.Range("K4").Value = $null
$ws1.Range("L4").Value = 2000000
$ws1.Range("M4").Value = "Đặng Ngọc Mai"
$ws1.Range("N4").Value = $null
$ws1.Range("O4").Value = 2000000
$ws1.Range("P4").Value = 0
$ws1.Range("Q4").Value = 2000000
$ws1.Range("R4").Value = 0
$ws1.Range("S4").Value = $null
$ws1.Range("T4").Value = $null
$ws1.Range("U4").Value = 0
$ws1.Range("V4").Value = 0

# Row 5: new data row.
Set-TextValue $ws1.Range("A5") "08-07-2024"
$ws1.Range("B5").Value = "HD-LUXURY"
$ws1.Range("C5").Value = 634
$ws1.Range("D5").Value = "LONG XUYÊN"
$ws1.Range("E5").Value = $null
$ws1.Range("F5").Value = $null
$ws1.Range("G5").Value = $null
$ws1.Range("H5").Value = $null
$ws1.Range("I5").Value = $null
$ws1.Range("J5").Value = $null
$ws1.Range("K5").Value = $null
$ws1.Range("L5").Value = 0
$ws1.Range("M5").Value = $null
$ws1.Range("N5").Value = $null
$ws1.Range("O5").Value = $null
$ws1.Range("P5").Value = 0
$ws1.Range("Q5").Value = 0
$ws1.Range("R5").Value = 0
$ws1.Range("S5").Value = $null
$ws1.Range("T5").Value = $null
$ws1.Range("U5").Value = $null
$ws1.Range("V5").Value = $null

# Row 6: new "Tổng" row.
$ws1.Range("A6").Value = $null
$ws1.Range("B6").Value = "Tổng"
$ws1.Range("C6").Value = 4
$ws1.Range("D6").Value = $null
$ws1.Range("E6").Value = $null
$ws1.Range("F6").Value = $null
$ws1.Range("G6").Value = $null
$ws1.Range("H6").Value = $null
$ws1.Range("I6").Value = 14000000
$ws1.Range("J6").Value = $null
$ws1.Range("K6").Value = 0
$ws1.Range("L6").Value = 14000000
$ws1.Range("M6").Value = $null
$ws1.Range("N6").Value = $null
$ws1.Range("O6").Value = 13000000
$ws1.Range("P6").Value = 0
$ws1.Range("Q6").Value = 13000000
$ws1.Range("R6").Value = 1000000
$ws1.Range("S6").Value = $null
$ws1.Range("T6").Value = $null
$ws1.Range("U6").Value = 100000
$ws1.Range("V6").Value = 0

# ============================================================
# Sheet 2: CHI TIẾT VỀ THU NỢ
# ============================================================
$ws2 = $wb.Worksheets.Item("CHI TIẾT VỀ THU NỢ")

# Row 2 used to be the "Tổng" row; it becomes a real data row.
$ws2.Range("A2").Value = "TN"
$ws2.Range("B2").Value = 174
Set-TextValue $ws2.Range("C2") "08-04-2024"
$ws2.Range("D2").Value = "LONG XUYÊN"
$ws2.Range("E2").Value = "HD-LUXURY-581"
Set-TextValue $ws2.Range("F2") "2024-07-23"
$ws2.Range("G2").Value = $null

# Row 3: new "Tổng" row.
$ws2.Range("A3").Value = "Tổng"
$ws2.Range("B3").Value = 1
$ws2.Range("C3").Value = $null
$ws2.Range("D3").Value = $null
$ws2.Range("E3").Value = $null
$ws2.Range("F3").Value = $null
$ws2.Range("G3").Value = 0

# ============================================================
# Sheet 3: CHI TIẾT CHI TIÊU
# ============================================================
$ws3 = $wb.Worksheets.Item("CHI TIẾT CHI TIÊU")

# New "Ghi chú" column header + blank column cells on the existing rows.
$ws3.Range("G1").Value = "Ghi chú"
$ws3.Range("G2").Value = $null
$ws3.Range("G3").Value = $null
$ws3.Range("G4").Value = $null
$ws3.Range("G5").Value = $null
$ws3.Range("G6").Value = $null

# Row 7 used to be the "Tổng" row; it becomes a real data row.
$ws3.Range("A7").Value = "CT"
$ws3.Range("B7").Value = 762
Set-TextValue $ws3.Range("C7") "08-06-2024"
$ws3.Range("D7").Value = "LONG XUYÊN"
$ws3.Range("E7").Value = "Ứng Lương"
$ws3.Range("F7").Value = 6000000
$ws3.Range("G7").Value = $null

# Row 8: new data row.
$ws3.Range("A8").Value = "CT"
$ws3.Range("B8").Value = 763
Set-TextValue $ws3.Range("C8") "08-06-2024"
$ws3.Range("D8").Value = "LONG XUYÊN"
$ws3.Range("E8").Value = "Chi Phí Sinh Hoạt Tại Cơ Sở"
$ws3.Range("F8").Value = 110000
$ws3.Range("G8").Value = $null

# Row 9: new data row.
$ws3.Range("A9").Value = "CT"
$ws3.Range("B9").Value = 765
Set-TextValue $ws3.Range("C9") "08-07-2024"
$ws3.Range("D9").Value = "LONG XUYÊN"
$ws3.Range("E9").Value = "Chi Phí CTV"
$ws3.Range("F9").Value = 2000000
$ws3.Range("G9").Value = $null

# Row 10: new "Tổng" row.
$ws3.Range("A10").Value = "Tổng"
$ws3.Range("B10").Value = 8
$ws3.Range("C10").Value = $null
$ws3.Range("D10").Value = $null
$ws3.Range("E10").Value = $null
$ws3.Range("F10").Value = 9967000
$ws3.Range("G10").Value = $null

# ============================================================
# Sheet 4: DOANH SỐ CÁ NHÂN
# ============================================================
$ws4 = $wb.Worksheets.Item("DOANH SỐ CÁ NHÂN")

$ws4.Range("A2").Value = "Nguyễn Hoàng Yến Quyên"
$ws4.Range("B2").Value = 0
$ws4.Range("C2").Value = 0
$ws4.Range("D2").Value = 11000000
$ws4.Range("E2").Value = 0
$ws4.Range("F2").Value = 0
$ws4.Range("G2").Value = 0
$ws4.Range("H2").Value = 0
$ws4.Range("I2").Value = 0
$ws4.Range("J2").Value = 0

$ws4.Range("A3").Value = "Nguyễn Phúc Nam"
$ws4.Range("B3").Value = 12000000
$ws4.Range("C3").Value = 0
$ws4.Range("D3").Value = 0
$ws4.Range("E3").Value = 0
$ws4.Range("F3").Value = 0
$ws4.Range("G3").Value = 0
$ws4.Range("H3").Value = 0
$ws4.Range("I3").Value = 0
$ws4.Range("J3").Value = 0

$ws4.Range("A4").Value = "Đào Vương Anh"
$ws4.Range("B4").Value = 0
$ws4.Range("C4").Value = 0
$ws4.Range("D4").Value = 0
$ws4.Range("E4").Value = 0
$ws4.Range("F4").Value = 2
$ws4.Range("G4").Value = 100000
$ws4.Range("H4").Value = 0
$ws4.Range("I4").Value = 0
$ws4.Range("J4").Value = 0

$ws4.Range("A5").Value = "Đặng Ngọc Mai"
$ws4.Range("B5").Value = 0
$ws4.Range("C5").Value = 0
$ws4.Range("D5").Value = 2000000
$ws4.Range("E5").Value = 0
$ws4.Range("F5").Value = 0
$ws4.Range("G5").Value = 0
$ws4.Range("H5").Value = 0
$ws4.Range("I5").Value = 0
$ws4.Range("J5").Value = 0

$ws4.Range("A6").Value = "Tổng"
$ws4.Range("B6").Value = 12000000
$ws4.Range("C6").Value = 0
$ws4.Range("D6").Value = 13000000
$ws4.Range("E6").Value = 0
$ws4.Range("F6").Value = 2
$ws4.Range("G6").Value = 100000
$ws4.Range("H6").Value = 0
$ws4.Range("I6").Value = 0
$ws4.Range("J6").Value = 0

# The old row 7 ("Tổng") is no longer needed now that row 6 holds it.
$ws4.Rows.Item(7).Delete()

# ============================================================
# Sheet 5: CHI TIÊU TỔNG HỢP
# ============================================================
$ws5 = $wb.Worksheets.Item("CHI TIÊU TỔNG HỢP")

$ws5.Range("B2").Value = 2500000
$ws5.Range("B3").Value = 610000
# Rows 4 (Chi Phí Vận Hành) and 5 (Trang thiết bị Y Tế) are unchanged.

# Row 6 used to be "Blank"; it becomes "Ứng Lương".
$ws5.Range("A6").Value = "Ứng Lương"
$ws5.Range("B6").Value = 6000000

# Row 7: new "Blank" row.
$ws5.Range("A7").Value = "Blank"
$ws5.Range("B7").Value = 0

# Row 8: "Tổng cộng" row (previously row 7).
$ws5.Range("A8").Value = "Tổng cộng"
$ws5.Range("B8").Value = 9967000

# ============================================================
# Sheet 6: LŨY KẾ NGÀY
# ============================================================
$ws6 = $wb.Worksheets.Item("LŨY KẾ NGÀY")

# Row 4 used to be the "Tổng" row; it becomes a real data row.
Set-TextValue $ws6.Range("A4") "08-04-2024"
$ws6.Range("B4").Value = 2000000
$ws6.Range("C4").Value = 2000000
$ws6.Range("D4").Value = 1
$ws6.Range("E4").Value = 0
$ws6.Range("F4").Value = 0
$ws6.Range("G4").Value = 2000000

# Row 5: new data row.
Set-TextValue $ws6.Range("A5") "08-06-2024"
$ws6.Range("B5").Value = 0
$ws6.Range("C5").Value = 0
$ws6.Range("D5").Value = 0
$ws6.Range("E5").Value = 0
$ws6.Range("F5").Value = 6110000
$ws6.Range("G5").Value = -6110000

# Row 6: new data row.
Set-TextValue $ws6.Range("A6") "08-07-2024"
$ws6.Range("B6").Value = 0
$ws6.Range("C6").Value = 0
$ws6.Range("D6").Value = 1
$ws6.Range("E6").Value = 0
$ws6.Range("F6").Value = 2000000
$ws6.Range("G6").Value = -2000000

# Row 7: new "Tổng" row.
$ws6.Range("A7").Value = "Tổng"
$ws6.Range("B7").Value = 14000000
$ws6.Range("C7").Value = 13000000
$ws6.Range("D7").Value = 4
$ws6.Range("E7").Value = 0
$ws6.Range("F7").Value = 9967000
$ws6.Range("G7").Value = 3033000

# ============================================================
# Sheet 7: QUỸ LƯƠNG
# ============================================================
$ws7 = $wb.Worksheets.Item("QUỸ LƯƠNG")

$ws7.Range("C4").Value = 857142.8571428573
$ws7.Range("C9").Value = 952380.9523809524
$ws7.Range("C10").Value = 428571.4285714285
$ws7.Range("C11").Value = 2857142.857142857
$ws7.Range("C12").Value = -2471428.571428571
$ws7.Range("C13").Value = 1517857.142857143
$ws7.Range("C14").Value = 1171428.571428572
$ws7.Range("C15").Value = 2122142.857142857
$ws7.Range("C16").Value = 995000
$ws7.Range("C22").Value = 9815952.380952382

# ============================================================
# Sheet 8: LỢI NHUẬN
# ============================================================
$ws8 = $wb.Worksheets.Item("LỢI NHUẬN")

$ws8.Range("B2").Value = 14000000
$ws8.Range("C2").Value = 13000000
$ws8.Range("D2").Value = 0.9285714285714286
$ws8.Range("E2").Value = 0.0714285714285714
$ws8.Range("F2").Value = 0
$ws8.Range("G2").Value = 13000000
$ws8.Range("H2").Value = 9967000
$ws8.Range("I2").Value = 9815952.380952382
$ws8.Range("J2").Value = 19782952.38095238
$ws8.Range("K2").Value = -6782952.380952381
$ws8.Range("L2").Value = -0.5217655677655677
